# Daily attendance processing - 2026-01-02 10:33:54
#
# Normalise the ordering of the "Recorded By" list in column G: each
# cell holds a comma-separated list of recorder identities (emails /
# "System" / "system"); reorder the entries of every such list into a
# fixed, stable priority order so the "real" recorder/account name
# consistently comes first, e.g.
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "backup@backdoor.com, System"         -> "System, backup@backdoor.com"
#   "admin@admin.com, System"             -> "System, admin@admin.com"
#   "system, backup@backdoor.com, System" -> "System, system, backup@backdoor.com"
# Lists that are already in the canonical order (or have a single
# entry) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 157
$col = 7   # column G ("Recorded By")

# Canonical priority order (lower = sorts earlier). Stable sort keeps
# relative order for any entries not listed here (none occur in this
# workbook, but this keeps the logic safe/generic).
$priorityOrder = @(
    "dnasr281@gmail.com",
    "System",
    "system",
    "admin@admin.com",
    "backup@backdoor.com"
)

function Get-Priority($name) {
    $idx = [Array]::IndexOf($priorityOrder, $name)
    if ($idx -lt 0) { return 999 }
    return $idx
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = @($text -split ", ")
    if ($parts.Count -le 1) {
        continue
    }

    # NOTE: Sort-Object with multiple -Property keys (e.g. "Pri, Orig")
    # does not reliably sort in this host, so combine priority + the
    # original index into a single composite numeric sort key instead
    # (priority dominates the high-order digits, original index breaks
    # ties and keeps the sort stable).
    $indexed = @()
    for ($i = 0; $i -lt $parts.Count; $i++) {
        $pri = Get-Priority $parts[$i]
        $key = ($pri * 1000) + $i
        $indexed += [PSCustomObject]@{ Value = $parts[$i]; Key = $key }
    }

    $sorted = $indexed | Sort-Object -Property Key

    $newParts = @()
    foreach ($item in $sorted) {
        $newParts += $item.Value
    }
    $newText = [string]::Join(", ", $newParts)

    if ($newText -ne $text) {
        $ws.Cells.Item($r, $col).Value = $newText
    }
}
